$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6500
$ws1.Range("F6").Value = 61
$ws1.Range("F9").Value = 92
$ws1.Range("F10").Value = 84
$ws1.Range("F15").Value = 3199
$ws1.Range("F17").Value = 201
$ws1.Range("F18").Value = 1866

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6500
$ws4.Range("F6").Value = 61
$ws4.Range("F10").Value = 92
$ws4.Range("F11").Value = 84
$ws4.Range("F16").Value = 3199
$ws4.Range("F18").Value = 201
$ws4.Range("F19").Value = 1866
